$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Tiny floating-point refresh of the existing row 3 date (matches source diff)
$ws.Range("A3").Value = 45804.4536994213

# New row 4 - same product, new price entry
$ws.Range("A4").Value = 45805.39372784663
$ws.Range("A4").NumberFormat = $ws.Range("A3").NumberFormat
$ws.Range("B4").Value = "EVOWHEY PROTEIN"
$ws.Range("C4").Value = "2Kg"
$ws.Range("D4").Value = "37,90€"
